# BIS-769: Fixed xls test files
# Add "Pattern" and "Pattern Type" header columns (M4, N4) to the sample
# property-type table, matching the style already used by the neighbouring
# "Unique" header cell (L4), and move the active selection onto the new
# columns (M4:N4) as in the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "Unique" header cell onto the two new
# header cells so they pick up the same style.
$ws.Range("L4").Copy() | Out-Null
$ws.Range("M4:N4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("M4").Value2 = "Pattern"
$ws.Range("N4").Value2 = "Pattern Type"

$ws.Range("M4:N4").Select() | Out-Null
